$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.683.86'
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").Value = '3.292.14'
$ws.Range("E3").Value = '  +0.59%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '185.44'
$ws.Range("E5").Value = '  +1.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '551.93'
$ws.Range("E6").Value = '  -0.71%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").Value = '3.286.81'
$ws.Range("E8").Value = '  +0.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.575'
$ws.Range("E9").Value = '  -2.99%  '

$ws.Range("E10").Value = '  -6.26%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.575'
$ws.Range("E11").Value = '  -2.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.56'
$ws.Range("E12").Value = '  -3.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000261'
$ws.Range("E13").Value = '  -1.69%  '

$ws.Range("D14").Value = '3.837.58'
$ws.Range("E14").Value = '  +1.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.41'
$ws.Range("E15").Value = '  -1.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '577.43'
$ws.Range("E16").Value = '  -8.88%  '

$ws.Range("D17").Value = '65.670.06'
$ws.Range("E17").Value = '  -0.03%  '

$ws.Range("E18").Value = '  +0.48%  '

$ws.Range("D19").Value = '3.312.32'
$ws.Range("E19").Value = '  +1.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.63'
$ws.Range("E20").Value = '  -1.27%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.79'
$ws.Range("E21").Value = '  -5.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.887'
$ws.Range("E22").Value = '  -1.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.96'
$ws.Range("E23").Value = '  +1.81%  '

$ws.Range("E24").Value = '  +0.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.87'
$ws.Range("E25").Value = '  -7.85%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.94'
$ws.Range("E26").Value = '  -0.91%  '

$ws.Range("E27").Value = '  +0.25%  '

$ws.Range("E28").Value = '  -2.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '30.29'
$ws.Range("E29").Value = '  -0.37%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.34'
$ws.Range("E30").Value = '  -3.82%  '

$ws.Range("E31").Value = '  +4.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '568.24'
$ws.Range("E32").Value = '  +4.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.67'
$ws.Range("E33").Value = '  -9.24%  '

$ws.Range("E34").Value = '  -1.91%  '

$ws.Range("E35").Value = '  -2.31%  '

$ws.Range("D36").Value = '3.720.01'
$ws.Range("E36").Value = '  +0.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.997'
$ws.Range("E37").Value = '  -0.36%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.19'
$ws.Range("E38").Value = '  -3.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '33.63'
$ws.Range("E39").Value = '  +3.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.125'
$ws.Range("E40").Value = '  -3.75%  '

$ws.Range("E41").Value = '  -5.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.11'
$ws.Range("E42").Value = '  -8.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.36'
$ws.Range("E43").Value = '  +3.69%  '

$ws.Range("E44").Value = '  -5.14%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.331'
$ws.Range("E45").Value = '  -1.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0405'
$ws.Range("E46").Value = '  -2.24%  '

$ws.Range("B47").Value = 'CoreDAO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.95'
$ws.Range("E47").Value = '  -11.94%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.127'
$ws.Range("E48").Value = '  -1.14%  '

$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.25%  '

$ws.Range("E50").Value = '  -4.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '127.16'
$ws.Range("E51").Value = '  +5.46%  '
